$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '21.733.04'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.32%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.540.00'
$ws.Range("D3").Style = "Normal"

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.001'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '290.15'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.22%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3894'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.50%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3183'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.68%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.11'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.61%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07198'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.47%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.060'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.74%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.01%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.640'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.25%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.65'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.40%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.630'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.96%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.540.94'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.87%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001104'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.01%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06578'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '83.18'
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.152'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.37%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.38'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.88%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.87'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.94%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.402'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +7.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '21.732.91'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.39%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.368'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.47%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '146.52'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.49%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.38'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.69%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.844'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.21%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.713.81'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.96%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '117.45'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.30%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.914'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.28%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9667'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -13.54%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08202'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.52%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.819'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.97%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06091'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.127'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.98%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02199'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.80%  '

$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2039'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.70%  '

$ws.Range("B40").Value = 'WEMIXTOKEN'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.450'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -12.91%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.187'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.57%  '

$ws.Range("E42").Value = '  -0.02%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.66'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.25%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5745'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.34%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.08'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.39%  '

$ws.Range("E46").Value = '  +0.40%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5489'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.36%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.159'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.37%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '116.40'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.05%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.875'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.75%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06719'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.96%  '
